# Update legacy GSC export data: roll the date window forward by 2 days.
#  - Drop the two oldest dates (2025-09-08, 2025-09-09) from the series.
#  - Append two new dates (2025-12-05, 2025-12-06) at the end.
#  - HTTPS URL counts (column C) shift up by 2 rows to line up with the
#    corresponding date; the two newly appended rows get a value of 0.
#  - Non-HTTPS URL counts (column B) remain 0 throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$lastRow = 89
$firstDataRow = 2

# Capture the existing column C (HTTPS URLs) values and the dates before
# overwriting anything, since the shift reads two rows ahead.
$oldC = @{}
$oldDates = @{}
for ($r = $firstDataRow; $r -le $lastRow; $r++) {
    $oldC[$r] = $ws.Cells.Item($r, 3).Value2
    $oldDates[$r] = [DateTime]::ParseExact($ws.Cells.Item($r, 1).Value2, "yyyy-MM-dd", $null)
}

# Shift dates forward by 2 days and shift the HTTPS URL counts to match,
# padding the two newest rows with 0 (no data collected yet).
for ($r = $firstDataRow; $r -le $lastRow; $r++) {
    $srcRow = $r + 2

    if ($srcRow -le $lastRow) {
        $newC = $oldC[$srcRow]
    } else {
        $newC = 0.0
    }

    $newDate = $oldDates[$r].AddDays(2)
    $dateText = $newDate.ToString("yyyy-MM-dd")

    # Dates are stored as plain text (General format), but Excel's normal
    # "typed" assignment auto-recognizes ISO date strings and converts them
    # to date serials. Force text entry by pre-formatting as Text, then
    # strip the formatting back off (ClearFormats leaves the stored value
    # untouched, only resets the cell's style back to the default/general
    # style) so the cell ends up identical in style to the rest of the
    # sheet while still holding literal text.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $dateText
    $cellA.ClearFormats()

    $ws.Cells.Item($r, 2).Value = 0.0
    $ws.Cells.Item($r, 3).Value = $newC
}
